$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "intervention_type" in K1, matching the style of the
# existing header row (A1:J1) by copying formats only.
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K1").Value = "intervention_type"

# New "intervention_type" values for rows 2-20 (column K)
$values = @(
    "PROCEDURE",
    "DEVICE",
    "PROCEDURE",
    "BEHAVIORAL",
    "PROCEDURE",
    "BEHAVIORAL",
    "DEVICE",
    "PROCEDURE",
    "OTHER",
    "PROCEDURE",
    "DIAGNOSTIC_TEST",
    "DIAGNOSTIC_TEST",
    "BEHAVIORAL",
    "BEHAVIORAL",
    "PROCEDURE",
    "OTHER",
    "DEVICE",
    "DRUG",
    "PROCEDURE"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $values[$i]
}
